$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.309.90"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "2.020.54"
$ws.Range("E3").Value = "  +6.46%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'244.36"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").Value = "'0.657"
$ws.Range("E6").Value = "  -5.02%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'44.15"
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("D9").Value = "'60.57"
$ws.Range("E9").Value = "  +6.44%  "
$ws.Range("D10").Value = "'0.357"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("E11").Value = "  -5.82%  "
$ws.Range("D12").Value = "'0.0981"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "'14.17"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").Value = "2.315.74"
$ws.Range("E14").Value = "  +6.44%  "
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "2.018.39"
$ws.Range("E16").Value = "  +6.33%  "
$ws.Range("E17").Value = "  -3.81%  "
$ws.Range("D18").Value = "36.377.66"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("D19").Value = "'70.76"
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("D21").Value = "'234.94"
$ws.Range("E21").Value = "  -4.90%  "
$ws.Range("D22").Value = "'12.61"
$ws.Range("E22").Value = "  -3.34%  "
$ws.Range("D23").Value = "'4.85"
$ws.Range("E23").Value = "  -6.58%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'2.43"
$ws.Range("E25").Value = "  -9.58%  "
$ws.Range("D26").Value = "'167.44"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "'8.62"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").Value = "'19.53"
$ws.Range("E28").Value = "  +6.19%  "
$ws.Range("D29").Value = "'1.92"
$ws.Range("E29").Value = "  -10.63%  "
$ws.Range("E30").Value = "  -6.51%  "
$ws.Range("D31").Value = "'21.55"
$ws.Range("E31").Value = "  +51.85%  "
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("D33").Value = "'0.0575"
$ws.Range("E33").Value = "  -5.01%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'1.87"
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("D36").Value = "'0.0858"
$ws.Range("E36").Value = "  +16.78%  "
$ws.Range("D37").Value = "'3.95"
$ws.Range("E37").Value = "  -7.45%  "
$ws.Range("E38").Value = "  +8.43%  "
$ws.Range("D39").Value = "'0.845"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("E40").Value = "  -11.87%  "
$ws.Range("E41").Value = "  -6.96%  "
$ws.Range("D42").Value = "'95.11"
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("D43").Value = "'1.10"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").Value = "'2.81"
$ws.Range("E44").Value = "  +17.05%  "
$ws.Range("D45").Value = "'15.50"
$ws.Range("E45").Value = "  -9.67%  "
$ws.Range("D46").Value = "1.301.41"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").Value = "'0.0815"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("D48").Value = "'2.78"
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D49").Value = "2.240.23"
$ws.Range("E49").Value = "  +7.95%  "
$ws.Range("E50").Value = "  -7.63%  "
$ws.Range("E51").Value = "  +14.62%  "
